# Refine API specs and other misc changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 becomes the "get one room" endpoint (rooms/create moves to row 11).
# Also fix D10's stray default style (s="0") so it matches the surrounding
# s="2" formatting -- a formats-only paste from a cell that already carries
# style 2 re-uses the existing style instead of minting a new one.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "/api/rooms/:id"
$ws.Range("C10").Value = "{data:{room:{}}}"
$ws.Range("C9").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 11: new "/api/rooms/create" endpoint row (rooms/:id now lives on row 10)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "/api/rooms/create"
$ws.Range("B11").Value = "{number,baseRent}"
$ws.Range("C11").Value = "{data:{room:{}}}"
$ws.Range("C10").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = "The newly created Room"

# ---------------------------------------------------------------------------
# Row 12: now a blank separator row (content moved up to row 11 / down to 13+)
# ---------------------------------------------------------------------------
$ws.Range("A12").ClearContents()
$ws.Range("C12").ClearContents()

# ---------------------------------------------------------------------------
# Row 13: /api/tenants (list)
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "/api/tenants"
$ws.Range("C13").Value = "{data:{tenants:[]}}"

# ---------------------------------------------------------------------------
# Row 14: /api/tenants/:id (get one)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "/api/tenants/:id"
$ws.Range("C14").Value = "{data:{tenant:{}}}"

# ---------------------------------------------------------------------------
# Row 15: /api/tenants/create
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "/api/tenants/create"
$ws.Range("B15").Value = "{name,phoneNumber,aadharCard,room}"
$ws.Range("C15").Value = "{data:{tenant:{}}}"
$ws.Range("D15").Value = "The newly created Tenant"

# ---------------------------------------------------------------------------
# Row 17: /api/transactions/?room (list) -- also fix stray D17 default style
# ---------------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("A17").Value = "/api/transactions/?room"
$ws.Range("C17").Value = "{data:{transactions:[]}}"

# ---------------------------------------------------------------------------
# Row 18: /api/transactions/:id (get one)
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "/api/transactions/:id"
$ws.Range("C18").Value = "{data:{transaction:{}}}"

# ---------------------------------------------------------------------------
# Row 19: /api/tenants/create/?roomNumber (create transaction)
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "/api/tenants/create/?roomNumber"
$ws.Range("B19").Value = "{room,balance,transfer,remarks}"
$ws.Range("C19").Value = "{data:{transaction:{}}}"
$ws.Range("D19").Value = "The newly created Transaction"

# ---------------------------------------------------------------------------
# Two new trailing blank rows (33-34), formatted like the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A33:AI34").Font.Name = "Times New Roman"
$ws.Range("A33:AI34").Font.Size = 14
$ws.Range("A33:AI34").HorizontalAlignment = 1

# ---------------------------------------------------------------------------
# Widen the first two columns to fit the longer endpoint/body text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 36.8
$ws.Columns.Item(2).ColumnWidth = 40.8

# ---------------------------------------------------------------------------
# Selection moved to B12 (the now-blank separator row).
# ---------------------------------------------------------------------------
$ws.Range("B12").Select()
